$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '300.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.67%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.43%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.024'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.91%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07949'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.96%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.904'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-10.98%'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.66%'

$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.034'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.83%'

$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.945'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '7.67%'

$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9218'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.50%'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1369'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '33.47%'

$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1846'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.32%'

$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09475'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.12%'

$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03599'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.14%'

$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09857'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.47%'

$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001392'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.06%'

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005739'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.16%'

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.512'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.02%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3425'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.73%'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.59%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.049'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.39%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2466'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '11.04%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04506'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.33%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.41%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004787'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.77%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001252'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.06%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003005'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-33.35%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01880'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.76%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04683'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.71%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007528'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-4.05%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009714'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '26.13%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1325'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.80%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002114'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.04%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009813'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-15.74%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006212'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.59%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.04%'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '64.88%'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-12.56%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.04%'
